$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style (date format) from BA1 to BB1
$ws.Range("BA1").Copy()
$ws.Range("BB1").PasteSpecial(-4122)

# Fill column BB (54) with values for rows 1-83
$ws.Cells.Item(1, 54).Value = 45986
$ws.Cells.Item(2, 54).Value = -0.3055163919209463
$ws.Cells.Item(3, 54).Value = -0.05930270183239372
$ws.Cells.Item(4, 54).Value = 0.6231390923824733
$ws.Cells.Item(5, 54).Value = -0.7175835844218028
$ws.Cells.Item(6, 54).Value = 0.6525336757609921
$ws.Cells.Item(7, 54).Value = 0.4125772253462401
$ws.Cells.Item(8, 54).Value = 0.0782632418494984
$ws.Cells.Item(9, 54).Value = 1.407630250560501
$ws.Cells.Item(10, 54).Value = -2.287736481071647
$ws.Cells.Item(11, 54).Value = 0.7460463752858573
$ws.Cells.Item(12, 54).Value = 0.4287266604711562
$ws.Cells.Item(13, 54).Value = -0.2522510312815029
$ws.Cells.Item(14, 54).Value = 0.8253228627186644
$ws.Cells.Item(15, 54).Value = -0.8956080111821194
$ws.Cells.Item(16, 54).Value = 0.6802092152583725
$ws.Cells.Item(17, 54).Value = -0.2702467235064887
$ws.Cells.Item(18, 54).Value = 0.2360788655137185
$ws.Cells.Item(19, 54).Value = 0.4710547017577227
$ws.Cells.Item(20, 54).Value = -0.6739669439274962
$ws.Cells.Item(21, 54).Value = -0.2065125907403171
$ws.Cells.Item(22, 54).Value = 0.1861693474447037
$ws.Cells.Item(23, 54).Value = 0.7530542497296722
$ws.Cells.Item(24, 54).Value = 0.3106180987640244
$ws.Cells.Item(25, 54).Value = 0.5709322285556908
$ws.Cells.Item(26, 54).Value = 1.247003875094862
$ws.Cells.Item(27, 54).Value = -0.776883278166693
$ws.Cells.Item(28, 54).Value = 1.384511819460911
$ws.Cells.Item(29, 54).Value = 0.09417825394891111
$ws.Cells.Item(30, 54).Value = -0.1288380179123294
$ws.Cells.Item(31, 54).Value = 0.1460429756367461
$ws.Cells.Item(32, 54).Value = -0.0194429241770564
$ws.Cells.Item(33, 54).Value = 0.2041957853813301
$ws.Cells.Item(34, 54).Value = 0.1772414135220401
$ws.Cells.Item(35, 54).Value = 0.7756459370471021
$ws.Cells.Item(36, 54).Value = 0.5676350012738425
$ws.Cells.Item(37, 54).Value = -0.7653316860800885
$ws.Cells.Item(38, 54).Value = 0.322227538137227
$ws.Cells.Item(39, 54).Value = -0.0579621925135001
$ws.Cells.Item(40, 54).Value = 0.6475935709367775
$ws.Cells.Item(41, 54).Value = 0.9603367340567104
$ws.Cells.Item(42, 54).Value = 0.3700548251239582
$ws.Cells.Item(43, 54).Value = 0.4443178743943008
$ws.Cells.Item(44, 54).Value = 0.5552956277764309
$ws.Cells.Item(45, 54).Value = 0.3556667645181193
$ws.Cells.Item(46, 54).Value = 0.6924895145077272
$ws.Cells.Item(47, 54).Value = 0.3519868436780342
$ws.Cells.Item(48, 54).Value = 0.489196258618918
$ws.Cells.Item(49, 54).Value = 0.1561519231779869
$ws.Cells.Item(50, 54).Value = 0.5297299217112936
$ws.Cells.Item(51, 54).Value = 0.7664191671286744
$ws.Cells.Item(52, 54).Value = 0.04527341468607915
$ws.Cells.Item(53, 54).Value = 0.06335359735614077
$ws.Cells.Item(54, 54).Value = 0.08598758370690973
$ws.Cells.Item(55, 54).Value = 0.2196343350075409
$ws.Cells.Item(56, 54).Value = 0.03810891122928695
$ws.Cells.Item(57, 54).Value = 0.4286160255108911
$ws.Cells.Item(58, 54).Value = 0.9
$ws.Cells.Item(59, 54).Value = 0.2
$ws.Cells.Item(60, 54).Value = 0.2
$ws.Cells.Item(61, 54).Value = 0
$ws.Cells.Item(62, 54).Value = -2.213339122522456
$ws.Cells.Item(63, 54).Value = -11.49785608241407
$ws.Cells.Item(64, 54).Value = 11.22930999924247
$ws.Cells.Item(65, 54).Value = -2.250986781122748
$ws.Cells.Item(66, 54).Value = -5.40098554941693
$ws.Cells.Item(67, 54).Value = 3.942549781810342
$ws.Cells.Item(68, 54).Value = 6.284666192508709
$ws.Cells.Item(69, 54).Value = -1.347757551663406
$ws.Cells.Item(70, 54).Value = 1.396500038188336
$ws.Cells.Item(71, 54).Value = 0.1885473380929312
$ws.Cells.Item(72, 54).Value = 1.224073604180177
$ws.Cells.Item(73, 54).Value = -1.674179157827311
$ws.Cells.Item(74, 54).Value = -0.6601113848982436
$ws.Cells.Item(75, 54).Value = 0.2658000717656392
$ws.Cells.Item(76, 54).Value = -0.01892863903084674
$ws.Cells.Item(77, 54).Value = 0.4261262404008619
$ws.Cells.Item(78, 54).Value = -0.07433428650158191
$ws.Cells.Item(79, 54).Value = 0.09298591595782568
$ws.Cells.Item(80, 54).Value = 0.2693878213604393
$ws.Cells.Item(81, 54).Value = 0.2408767182737677
$ws.Cells.Item(82, 54).Value = 0.5597354586130052
$ws.Cells.Item(83, 54).Value = 0.1368731201391853

# Add new row 84: A84 date, BB84 value; copy A-column date style from A83
$ws.Range("A83").Copy()
$ws.Range("A84").PasteSpecial(-4122)
$ws.Cells.Item(84, 1).Value = 45884
$ws.Cells.Item(84, 54).Value = -0.2551464291630765

$excel.CutCopyMode = 0
